$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (date "Förändrad") for rows 2-5: 45243 -> 45244 (increment by one day)
$ws.Range("C2:C5").Value = 45244
